$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

$ws.Columns.Item(5).ColumnWidth = 13

$ws.Range("F3").Select() | Out-Null
